$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 5: new entry - 21.03.2021, 0.5h, "Fixed git-merging issues"
$ws.Cells.Item(5, 1).Value = (Get-Date -Year 2021 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(5, 2).Value = 0.5
$ws.Cells.Item(5, 4).Value = "Fixed git-merging issues"

# Row 6: new entry - 05.04.2021, 6h, "Organizational stuff; Research on jquery-ui"
$ws.Cells.Item(6, 1).Value = (Get-Date -Year 2021 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(6, 2).Value = 6
$ws.Cells.Item(6, 4).Value = "Organizational stuff; Research on jquery-ui"

# Update the selected cell/range shown in the sheet view
$ws.Range("H6").Select()

$wb.Save()
